# Generate Report for Handback
# ------------------------------------------------------------------
# This applies the "handback" report update to localization-status.xlsx:
#   * Status text updates from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shared everywhere it's used).
#   * The zh-cn and de-de sheets gain two new populated columns:
#       F = "Latest Target File"   (the source .md file, same link as col A)
#       G = "Latest Handback File" (the generated .xlf file, same link as col D)
#   * The "Latest Handback DateTime" (col H), previously the zero-date
#     placeholder, is stamped with the real handback time - zh-cn finished
#     a few seconds before de-de.

$wb = $excel.ActiveWorkbook

$mdUrlBase  = "https://github.com/OpenLocalizationTest/oltest/blob/565e415293c24513d52119f35364207e71fc5974/e2e"
$zhXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fbd4caba024e55aff767d327b8a6c251a557049/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high"
$deXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/31c1c19b075ddc4b3b19846270328bb13461b302/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high"

$file1 = "54855bbe-89c5-45c0-b19d-0387efbfd8bc"
$file2 = "653d20ae-e458-478b-8cf9-8fbbac8cbb49"

$zhXlf1 = "$file1.512776ff87e1407d88d21f8417fcb5064fdeeda1.zh-cn.xlf"
$zhXlf2 = "$file2.6fcd863d1d3ad255908559c1693faf9e72c97409.zh-cn.xlf"
$deXlf1 = "$file1.512776ff87e1407d88d21f8417fcb5064fdeeda1.de-de.xlf"
$deXlf2 = "$file2.6fcd863d1d3ad255908559c1693faf9e72c97409.de-de.xlf"

# ---- 1. Status wording: "Ready for handoff" -> "Handed back: in sync with en-US"
# The same shared string is used on every sheet's Status column, so a
# straightforward Find/Replace across the whole workbook reproduces the
# sharedStrings.xml text edit everywhere it's referenced.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $used.Replace("Ready for handoff", "Handed back: in sync with en-US", 1, 1, $false, $false, $false) | Out-Null
}

# ---- 2. zh-cn sheet: fill in the "Latest Target File" / "Latest Handback File"
#         columns (F/G) for both data rows, and stamp the handback time.
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("F2").Value = "$file1.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "$mdUrlBase/$file1.md", "", "", "$file1.md") | Out-Null
$ws.Range("F2").Style = $ws.Range("A2").Style

$ws.Range("G2").Value = $zhXlf1
$ws.Hyperlinks.Add($ws.Range("G2"), "$zhXlfBase/$zhXlf1", "", "", $zhXlf1) | Out-Null
$ws.Range("G2").Style = $ws.Range("D2").Style

$ws.Range("F3").Value = "$file2.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "$mdUrlBase/$file2.md", "", "", "$file2.md") | Out-Null
$ws.Range("F3").Style = $ws.Range("A3").Style

$ws.Range("G3").Value = $zhXlf2
$ws.Hyperlinks.Add($ws.Range("G3"), "$zhXlfBase/$zhXlf2", "", "", $zhXlf2) | Out-Null
$ws.Range("G3").Style = $ws.Range("D3").Style

$ws.Range("H2").Value = "2016-03-12 14:11:38"
$ws.Range("H3").Value = "2016-03-12 14:11:38"

# ---- 3. de-de sheet: same treatment, slightly later handback timestamp.
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("F2").Value = "$file1.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "$mdUrlBase/$file1.md", "", "", "$file1.md") | Out-Null
$ws.Range("F2").Style = $ws.Range("A2").Style

$ws.Range("G2").Value = $deXlf1
$ws.Hyperlinks.Add($ws.Range("G2"), "$deXlfBase/$deXlf1", "", "", $deXlf1) | Out-Null
$ws.Range("G2").Style = $ws.Range("D2").Style

$ws.Range("F3").Value = "$file2.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "$mdUrlBase/$file2.md", "", "", "$file2.md") | Out-Null
$ws.Range("F3").Style = $ws.Range("A3").Style

$ws.Range("G3").Value = $deXlf2
$ws.Hyperlinks.Add($ws.Range("G3"), "$deXlfBase/$deXlf2", "", "", $deXlf2) | Out-Null
$ws.Range("G3").Style = $ws.Range("D3").Style

$ws.Range("H2").Value = "2016-03-12 14:11:44"
$ws.Range("H3").Value = "2016-03-12 14:11:44"

Write-Output "Handback report generated"
